$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B11").Value = 'https://media.boohoo.com/i/boohoo/bmm71158_charcoal_xl?w=900&qlt=default&fmt.jp2.qlt=70&fmt=auto&sm=fit'
$ws.Range("D11").Value = 'REGULAR CREW NECK TWO TONE RIB KNITTED JUMPER'
$ws.Range("E11").Value = '$30.00'
$ws.Range("F11").Value = 'https://ca.boohoo.com/regular-crew-neck-two-tone-rib-knitted-jumper/BMM71158.html?color=115'
$ws.Range("B12").Value = 'https://media.boohoo.com/i/boohoo/bmm76020_grey_xl?w=900&qlt=default&fmt.jp2.qlt=70&fmt=auto&sm=fit'
$ws.Range("D12").Value = 'PLUS CORE COLOUR BLOCK SWEATER TRACKSUIT'
$ws.Range("E12").Value = '$66.00'
$ws.Range("F12").Value = 'https://ca.boohoo.com/plus-core-colour-block-sweater-tracksuit/BMM76020.html?color=131'
$ws.Range("B13").Value = 'https://media.boohoo.com/i/boohoo/bmm65956_stone_xl?w=900&qlt=default&fmt.jp2.qlt=70&fmt=auto&sm=fit'
$ws.Range("D13").Value = 'SHORT SLEEVE MUSCLE FIT STRIPE KNIT POLO'
$ws.Range("E13").Value = '$36.50'
$ws.Range("F13").Value = 'https://ca.boohoo.com/short-sleeve-muscle-fit-stripe-knit-polo/BMM65956.html?color=165'
$ws.Range("B14").Value = 'https://media.boohoo.com/i/boohoo/bmm58525_grey_xl?w=900&qlt=default&fmt.jp2.qlt=70&fmt=auto&sm=fit'
$ws.Range("D14").Value = 'BOXY DROP SHOULDER EYE GRAPHIC KNITTED JUMPER'
$ws.Range("E14").Value = '$15.00'
$ws.Range("F14").Value = 'https://ca.boohoo.com/boxy-drop-shoulder-eye-graphic-knitted-jumper/BMM58525.html'
$ws.Range("B15").Value = 'https://media.boohoo.com/i/boohoo/bmm74381_charcoal_xl?w=900&qlt=default&fmt.jp2.qlt=70&fmt=auto&sm=fit'
$ws.Range("D15").Value = 'PLUS TEXTURED KNITTED POLO IN CHARCOAL'
$ws.Range("E15").Value = '$40.00'
$ws.Range("F15").Value = 'https://ca.boohoo.com/plus-textured-knitted-polo-in-charcoal/BMM74381.html?color=115'
$ws.Range("B16").Value = 'https://media.boohoo.com/i/boohoo/bmm49664_navy_xl?w=900&qlt=default&fmt.jp2.qlt=70&fmt=auto&sm=fit'
$ws.Range("D16").Value = 'V NECK STRIPED KNITTED POLO'
$ws.Range("E16").Value = '$50.00'
$ws.Range("F16").Value = 'https://ca.boohoo.com/v-neck-striped-knitted-polo/BMM49664.html'
$ws.Range("B17").Value = 'https://media.boohoo.com/i/boohoo/bmm68577_black_xl?w=900&qlt=default&fmt.jp2.qlt=70&fmt=auto&sm=fit'
$ws.Range("D17").Value = 'BOXY FIT KNITTED MOTO SHIRT'
$ws.Range("E17").Value = '$36.50'
$ws.Range("F17").Value = 'https://ca.boohoo.com/boxy-fit-knitted-moto-shirt/BMM68577.html?color=105'
$ws.Range("B18").Value = 'https://media.boohoo.com/i/boohoo/bmm75278_ecru_xl?w=900&qlt=default&fmt.jp2.qlt=70&fmt=auto&sm=fit'
$ws.Range("D18").Value = 'LONG SLEEVE COLOUR BLOCK KNITTED SHIRT'
$ws.Range("E18").Value = '$40.00'
$ws.Range("F18").Value = 'https://ca.boohoo.com/long-sleeve-colour-block-knitted-shirt/BMM75278.html?color=124'
$ws.Range("B19").Value = 'https://media.boohoo.com/i/boohoo/bmm58824_ecru_xl?w=900&qlt=default&fmt.jp2.qlt=70&fmt=auto&sm=fit'
$ws.Range("D19").Value = 'BOXY KNITTED RIBBED HOODIE'
$ws.Range("E19").Value = '$46.50'
$ws.Range("F19").Value = 'https://ca.boohoo.com/boxy-knitted-ribbed-hoodie/BMM58824.html?color=124'
$ws.Range("B20").Value = 'https://media.boohoo.com/i/boohoo/bmm74962_ecru_xl?w=900&qlt=default&fmt.jp2.qlt=70&fmt=auto&sm=fit'
$ws.Range("D20").Value = 'SHORT SLEEVE BOXY FIT REVERE OPEN KNIT POLO IN ECRU'
$ws.Range("E20").Value = '$40.00'
$ws.Range("F20").Value = 'https://ca.boohoo.com/short-sleeve-boxy-fit-revere-open-knit-polo-in-ecru/BMM74962.html'
$ws.Range("B21").Value = 'https://media.boohoo.com/i/boohoo/bmm61582_stone_xl?w=900&qlt=default&fmt.jp2.qlt=70&fmt=auto&sm=fit'
$ws.Range("D21").Value = 'OVERSIZED HOMME BRUSHED RIB KNIT VEST'
$ws.Range("E21").Value = '$37.00'
$ws.Range("F21").Value = 'https://ca.boohoo.com/oversized-homme-brushed-rib-knit-vest-/BMM61582.html'
$ws.Range("B22").Value = 'https://media.boohoo.com/i/boohoo/bmm80656_black_xl?w=900&qlt=default&fmt.jp2.qlt=70&fmt=auto&sm=fit'
$ws.Range("D22").Value = 'OVERSIZED BOXY CROCHET KNIT POLO IN BLACK'
$ws.Range("E22").Value = '$40.00'
$ws.Range("F22").Value = 'https://ca.boohoo.com/oversized-boxy-crochet-knit-polo-in-black/BMM80656.html?color=105'
$ws.Range("B23").Value = 'https://media.boohoo.com/i/boohoo/bmm72243_chocolate_xl?w=900&qlt=default&fmt.jp2.qlt=70&fmt=auto&sm=fit'
$ws.Range("D23").Value = 'FLUFFY KNITTED BOXY HOODIE'
$ws.Range("E23").Value = '$27.00'
$ws.Range("F23").Value = 'https://ca.boohoo.com/fluffy-knitted-boxy-hoodie/BMM72243.html?color=186'
$ws.Range("B24").Value = 'https://media.boohoo.com/i/boohoo/bmm74797_slate%20grey_xl?w=900&qlt=default&fmt.jp2.qlt=70&fmt=auto&sm=fit'
$ws.Range("D24").Value = 'REGULAR FIT WAFFLE KNIT JUMPER'
$ws.Range("E24").Value = '$36.50'
$ws.Range("F24").Value = 'https://ca.boohoo.com/regular-fit-waffle-knit-jumper/BMM74797.html?color=849'
$ws.Range("B25").Value = 'https://media.boohoo.com/i/boohoo/bmm69862_ecru_xl?w=900&qlt=default&fmt.jp2.qlt=70&fmt=auto&sm=fit'
$ws.Range("D25").Value = 'RELAXED CABLE KNITTED VEST'
$ws.Range("E25").Value = '$37.00'
$ws.Range("F25").Value = 'https://ca.boohoo.com/relaxed-cable-knitted-vest-/BMM69862.html?color=124'
$ws.Range("B26").Value = 'https://media.boohoo.com/i/boohoo/bmm69872_chocolate_xl?w=900&qlt=default&fmt.jp2.qlt=70&fmt=auto&sm=fit'
$ws.Range("D26").Value = 'MUSCLE FIT ZIP THROUGH RIB KNIT JACKET'
$ws.Range("E26").Value = '$37.00'
$ws.Range("F26").Value = 'https://ca.boohoo.com/muscle-fit-zip-through-rib-knit-jacket/BMM69872.html?color=186'
$ws.Range("B27").Value = 'https://media.boohoo.com/i/boohoo/bmm68476_ecru_xl?w=900&qlt=default&fmt.jp2.qlt=70&fmt=auto&sm=fit'
$ws.Range("D27").Value = 'BOXY FLUFFY KNITTED CARDIGAN'
$ws.Range("E27").Value = '$60.00'
$ws.Range("F27").Value = 'https://ca.boohoo.com/boxy-fluffy-knitted-cardigan/BMM68476.html?color=124'
$ws.Range("B28").Value = 'https://media.boohoo.com/i/boohoo/bmm59722_ecru_xl?w=900&qlt=default&fmt.jp2.qlt=70&fmt=auto&sm=fit'
$ws.Range("D28").Value = 'OVERSIZED CREW NECK FLUFFY KNITTED JUMPER'
$ws.Range("E28").Value = '$50.00'
$ws.Range("F28").Value = 'https://ca.boohoo.com/oversized-crew-neck-fluffy-knitted-jumper/BMM59722.html?color=124'
